$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 6; this shifts the existing rows 6..122 down to 7..123
# and keeps the dimension growing from A1:R122 to A1:R123.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with a new weekly record. Columns that are
# not explicitly part of the edit (A, B, C, E, F, G, H, I, K, L, M, N, O, Q, R)
# keep the same values as the row directly below (the former row 6, now row 7),
# while D (Fecha) and J (Volumen) take on the new values from the commit.
$ws.Cells.Item(6, 1).Value = 2
$ws.Cells.Item(6, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(6, 3).Value = "Coquimbo"
$ws.Cells.Item(6, 4).Value = 45245
$ws.Cells.Item(6, 5).Value = 4
$ws.Cells.Item(6, 6).Value = 100112026
$ws.Cells.Item(6, 7).Value = "Haba"
$ws.Cells.Item(6, 8).Value = "Sin especificar"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 500
$ws.Cells.Item(6, 11).Value = 8000
$ws.Cells.Item(6, 12).Value = 9000
$ws.Cells.Item(6, 13).Value = 8500
$ws.Cells.Item(6, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(6, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(6, 16).Value = 340
$ws.Cells.Item(6, 17).Value = 25
$ws.Cells.Item(6, 18).Value = "Hortaliza"
